$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column price cells are treated as text (preserve exact string formatting)
# by pre-setting the number format to Text for the whole price column range, then
# reverting the style to Normal afterwards so no extra styling is introduced.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "44.555.44"
$ws.Range("E2").Value = "  +3.70%  "
$ws.Range("D3").Value = "2.424.00"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "312.55"
$ws.Range("E5").Value = "  +3.34%  "
$ws.Range("D6").Value = "101.66"
$ws.Range("E6").Value = "  +6.02%  "
$ws.Range("D7").Value = "0.512"
$ws.Range("E7").Value = "  +1.79%  "
$ws.Range("E9").Value = "  +4.89%  "
$ws.Range("E10").Value = "  +3.40%  "
$ws.Range("D11").Value = "0.0801"
$ws.Range("E11").Value = "  +1.99%  "
$ws.Range("E12").Value = "  +1.75%  "
$ws.Range("D13").Value = "18.91"
$ws.Range("E13").Value = "  +2.91%  "
$ws.Range("E14").Value = "  +2.90%  "
$ws.Range("D15").Value = "2.803.88"
$ws.Range("D16").Value = "2.422.38"
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("E17").Value = "  +4.85%  "
$ws.Range("D18").Value = "44.454.03"
$ws.Range("E18").Value = "  +3.53%  "
$ws.Range("D19").Value = "12.50"
$ws.Range("E19").Value = "  +3.99%  "
$ws.Range("E20").Value = "  +2.11%  "
$ws.Range("D21").Value = "0.0₃0924"
$ws.Range("E21").Value = "  +4.40%  "
$ws.Range("D22").Value = "68.87"
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("D23").Value = "241.01"
$ws.Range("E24").Value = "  +4.90%  "
$ws.Range("D25").Value = "2.46"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("D27").Value = "25.17"
$ws.Range("E27").Value = "  +2.58%  "
$ws.Range("E28").Value = "  -4.34%  "
$ws.Range("D29").Value = "9.61"
$ws.Range("E29").Value = "  +3.60%  "
$ws.Range("D30").Value = "33.30"
$ws.Range("E30").Value = "  +5.22%  "
$ws.Range("D31").Value = "48.60"
$ws.Range("E31").Value = "  +1.29%  "
$ws.Range("E32").Value = "  +18.41%  "
$ws.Range("D33").Value = "19.57"
$ws.Range("E34").Value = "  +3.08%  "
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("D36").Value = "0.0767"
$ws.Range("E36").Value = "  +6.94%  "
$ws.Range("E37").Value = "  +4.17%  "
$ws.Range("E38").Value = "  +2.77%  "
$ws.Range("E39").Value = "  +4.15%  "
$ws.Range("D40").Value = "124.18"
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("D42").Value = "2.19"
$ws.Range("E42").Value = "  -4.22%  "
$ws.Range("D43").Value = "21.34"
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").Value = "0.0290"
$ws.Range("E44").Value = "  +4.01%  "
$ws.Range("D45").Value = "1.950.24"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "2.93"
$ws.Range("E46").Value = "  +7.76%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "9.52"
$ws.Range("E47").Value = "  +3.59%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "1.65"
$ws.Range("E48").Value = "  +9.42%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "54.35"
$ws.Range("E49").Value = "  +5.23%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "73.93"
$ws.Range("E50").Value = "  +3.59%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "4.64"
$ws.Range("E51").Value = "  +5.78%  "

# Revert the price column style back to the default ("Normal") now that the
# values have been written as text, so no extra cell styling remains applied.
$ws.Range("D2:D51").Style = "Normal"

